$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (columns B through Q) - identical for every data row (2-26)
$newValues = @(
    [double]"0.9999677858346675",
    [double]"0.9990493920370204",
    [double]"0.9999999999998502",
    [double]"0.999933657405726",
    [double]"0.9999661959859565",
    [double]"3.007050929946918e-05",
    [double]"0.0008873508065744985",
    [double]"8.629902877618854e-14",
    [double]"7.105010605329052e-05",
    [double]"3.552505294840163e-05",
    [double]"0.0003164627591308915",
    [double]"0.005483658386466938",
    [double]"1.000026659998896",
    [double]"0.005717109076716589",
    [double]"126.8239312479565",
    [double]"191.4243499659711"
)

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $newValues.Length; $i++) {
        $col = 2 + $i  # column B is index 2
        $ws.Cells.Item($row, $col).Value = $newValues[$i]
    }
}
